$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7531644701957703
$ws.Range("B1").Value = 1.412539720535278
$ws.Range("C1").Value = 5.331111431121826
$ws.Range("D1").Value = 3.176127910614014
$ws.Range("E1").Value = 1.519210696220398
